# Finishing the hyperparameter and formula tuning vignette:
# restructure the "workflow" control sheet to add saved-hyperparameter,
# seasonality-search and random-effect-search rows, and drop the old
# "add_trend" row (and its yellow highlight style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workflow")

# Remove the old "add_trend" row entirely -- this also drops the yellow
# highlight fill/style that only that row used.
$ws.Rows.Item(6).Delete()

# Clear out the remaining old rows (2:7) so we can lay down the new,
# reordered/expanded set of rows from scratch.
$ws.Range("A2:C7").ClearContents()

# Row 2: which variable are we modeling (moved up from old row 5) -- all
# of these strings already exist in the shared-string table.
$ws.Range("A2").Value = "Which variable are we modeling?"
$ws.Range("B2").Value = "sales"
$ws.Range("C2").Value = "Y"

# Row 3: run hyper parameter tuning process (moved down from old row 2)
$ws.Range("A3").Value = "Run a hyper parameter tuning process?  If FALSE, will use saved hyper-parameters"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = "tune_this_time"

# New rows/strings below -- write them in the same order the saved
# workbook introduces them into the shared-string table.
$ws.Range("A4").Value = "saved hyper parameter tibble"
$ws.Range("C4").Value = "saved_hypers_filename"

$ws.Range("A5").Value = "Find best seasonality spec?"
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = "search_seasonality"

$ws.Range("A6").Value = "seasonality interactions to test?"
$ws.Range("B6").Value = "store"
$ws.Range("C6").Value = "interaction_fft"

# Row 7: number of fourier terms (moved from old row 3, value 0 -> 1)
$ws.Range("A7").Value = "number of fourier terms to use for seasonality (up to 5; 2 is roughly semi-annual + trimesters, 3 is semi_annual, trimester,quarterly)"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "fft_terms"

# Row 8: should different random effects be tested (new). C8 ("search_randoms")
# and B10 ("(TV1|store)") land in the shared-string table before A8's text,
# so write them first to match.
$ws.Range("C8").Value = "search_randoms"
$ws.Range("B10").Value = "(TV1|store)"
$ws.Range("A8").Value = "should different random effects be tested?"
$ws.Range("B8").Value = $true

# Row 9: random effects intercepts (moved from old row 4)
$ws.Range("A9").Value = "random effects intercepts -- comma separated list of categorical vars that should have a random effect"
$ws.Range("B9").Value = "(1|store)"
$ws.Range("C9").Value = "list_rand_ints"

# Row 10: random slopes (moved from old row 7, value trimmed to single term)
$ws.Range("A10").Value = "random slopes"
$ws.Range("C10").Value = "list_rand_slopes"

# Update selection to match the saved workbook (active cell A8)
$ws.Range("A8").Select()
